$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Fix the "earier" typo -> "earlier" in the EventID paragraph. This
#    also clears the (now stale) spell-check proofErr markers that used
#    to flag the misspelling.
# ---------------------------------------------------------------------
$r1 = $d.Content
$r1.Find.Execute("is a UUID generated in an earier step", $true, $false, $false, `
                  $false, $false, $true, 1, $false, `
                  "is a UUID generated in an earlier step", 2) | Out-Null

# ---------------------------------------------------------------------
# 2. Word drops a "_GoBack" bookmark at the most recent edit location.
#    Move it from its old spot (near the end of the document) onto the
#    empty paragraph immediately below the paragraph we just edited.
# ---------------------------------------------------------------------
$goBackTarget = $d.Paragraphs(7).Range
$d.Bookmarks.Add("_GoBack", $goBackTarget)

# ---------------------------------------------------------------------
# 3. Near the end of the document, "... locationIDs and the[bookmark]n
#    it's ready for upsert." becomes "... locationIDs and then it's
#    ready for upsert." -- i.e. "the" + "n" joins into "then", merging
#    those two runs and dropping the now-superseded bookmark that used
#    to sit between them.
# ---------------------------------------------------------------------
$r3 = $d.Content
$r3.Find.Execute("and then it" + [char]0x2019 + "s ready for", $true, $false, $false, `
                  $false, $false, $true, 1, $false, `
                  "and then it" + [char]0x2019 + "s ready for", 2) | Out-Null
